$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Cells.Item(2, 7).Value = 2.98
$ws.Cells.Item(2, 8).Value = 2.68
$ws.Cells.Item(2, 10).Value = 3.15
$ws.Cells.Item(2, 11).Value = 3.65
$ws.Cells.Item(2, 23).Value = 1.5
$ws.Cells.Item(2, 28).Value = 12.5

# Row 3
$ws.Cells.Item(3, 10).Value = 1.03
$ws.Cells.Item(3, 14).Value = 1.11
$ws.Cells.Item(3, 20).Value = 1.03
$ws.Cells.Item(3, 21).Value = 1.03

# Row 4
$ws.Cells.Item(4, 6).Value = 1.78
$ws.Cells.Item(4, 7).Value = 1.93
$ws.Cells.Item(4, 9).Value = 5.7
$ws.Cells.Item(4, 11).Value = 3.95
$ws.Cells.Item(4, 14).Value = 3.15
$ws.Cells.Item(4, 15).Value = 1.37
$ws.Cells.Item(4, 16).Value = 1.74
$ws.Cells.Item(4, 17).Value = 2.1
$ws.Cells.Item(4, 18).Value = 1.28
$ws.Cells.Item(4, 19).Value = 3.9
$ws.Cells.Item(4, 21).Value = 1.86
$ws.Cells.Item(4, 22).Value = 1.21
$ws.Cells.Item(4, 24).Value = 14.5
$ws.Cells.Item(4, 25).Value = 19
$ws.Cells.Item(4, 29).Value = 9.800000000000001
$ws.Cells.Item(4, 30).Value = 23
$ws.Cells.Item(4, 31).Value = 100
$ws.Cells.Item(4, 32).Value = 11.5
$ws.Cells.Item(4, 33).Value = 11
$ws.Cells.Item(4, 34).Value = 23
$ws.Cells.Item(4, 40).Value = 18

# Row 5
$ws.Cells.Item(5, 7).Value = 1.64
$ws.Cells.Item(5, 9).Value = 7.2
$ws.Cells.Item(5, 10).Value = 4.1
$ws.Cells.Item(5, 19).Value = 3.15
$ws.Cells.Item(5, 23).Value = 2.56
$ws.Cells.Item(5, 25).Value = 980
$ws.Cells.Item(5, 26).Value = 65
$ws.Cells.Item(5, 30).Value = 980
$ws.Cells.Item(5, 37).Value = 980
$ws.Cells.Item(5, 38).Value = 980

# Row 6
$ws.Cells.Item(6, 6).Value = 2.06
$ws.Cells.Item(6, 7).Value = 2.74
$ws.Cells.Item(6, 8).Value = 2.9
$ws.Cells.Item(6, 9).Value = 5.4
$ws.Cells.Item(6, 11).Value = 6
$ws.Cells.Item(6, 22).Value = 1.22
$ws.Cells.Item(6, 23).Value = 1.58

# Row 7
$ws.Cells.Item(7, 6).Value = 5.7
$ws.Cells.Item(7, 7).Value = 6.8
$ws.Cells.Item(7, 9).Value = 1.79
$ws.Cells.Item(7, 22).Value = 2.26
$ws.Cells.Item(7, 23).Value = 1.18
$ws.Cells.Item(7, 25).Value = 8.800000000000001
$ws.Cells.Item(7, 27).Value = 22
$ws.Cells.Item(7, 28).Value = 22
$ws.Cells.Item(7, 31).Value = 25

# Row 9
$ws.Cells.Item(9, 6).Value = 1.59
$ws.Cells.Item(9, 7).Value = 1.92
$ws.Cells.Item(9, 9).Value = 9.4
$ws.Cells.Item(9, 10).Value = 3.1
$ws.Cells.Item(9, 11).Value = 5.1
$ws.Cells.Item(9, 17).Value = 2.44
$ws.Cells.Item(9, 19).Value = 2.44
$ws.Cells.Item(9, 22).Value = 1.12
$ws.Cells.Item(9, 23).Value = 2.08

# Row 10
$ws.Cells.Item(10, 8).Value = 2.4
$ws.Cells.Item(10, 9).Value = 3.35
$ws.Cells.Item(10, 10).Value = 2.84
$ws.Cells.Item(10, 22).Value = 1.48
$ws.Cells.Item(10, 23).Value = 1.44
$ws.Cells.Item(10, 25).Value = 1000

# Row 11
$ws.Cells.Item(11, 7).Value = 2.3
$ws.Cells.Item(11, 9).Value = 4.9
$ws.Cells.Item(11, 10).Value = 2.84
$ws.Cells.Item(11, 14).Value = 2.2
$ws.Cells.Item(11, 15).Value = 1.78
$ws.Cells.Item(11, 16).Value = 1.37
$ws.Cells.Item(11, 17).Value = 3.35
$ws.Cells.Item(11, 19).Value = 8.199999999999999
$ws.Cells.Item(11, 20).Value = 2.6
$ws.Cells.Item(11, 23).Value = 1.77
$ws.Cells.Item(11, 24).Value = 6.4
$ws.Cells.Item(11, 25).Value = 10
$ws.Cells.Item(11, 27).Value = 150
$ws.Cells.Item(11, 28).Value = 5.9
$ws.Cells.Item(11, 29).Value = 7.4
$ws.Cells.Item(11, 30).Value = 23
$ws.Cells.Item(11, 37).Value = 42
$ws.Cells.Item(11, 38).Value = 1000
$ws.Cells.Item(11, 39).Value = 400
$ws.Cells.Item(11, 40).Value = 48

# Row 12
$ws.Cells.Item(12, 7).Value = 2.94
$ws.Cells.Item(12, 8).Value = 2.88
$ws.Cells.Item(12, 9).Value = 3.2
$ws.Cells.Item(12, 13).Value = 1.11
$ws.Cells.Item(12, 22).Value = 1.45
$ws.Cells.Item(12, 32).Value = 20
$ws.Cells.Item(12, 34).Value = 25
$ws.Cells.Item(12, 36).Value = 55
$ws.Cells.Item(12, 37).Value = 44

# Row 13
$ws.Cells.Item(13, 6).Value = 2.74
$ws.Cells.Item(13, 7).Value = 2.98
$ws.Cells.Item(13, 9).Value = 3
$ws.Cells.Item(13, 10).Value = 3.15
$ws.Cells.Item(13, 13).Value = 1.08
$ws.Cells.Item(13, 19).Value = 3.9
$ws.Cells.Item(13, 22).Value = 1.5
$ws.Cells.Item(13, 24).Value = 14.5
$ws.Cells.Item(13, 35).Value = 60
$ws.Cells.Item(13, 38).Value = 60
$ws.Cells.Item(13, 41).Value = 42

# Row 14
$ws.Cells.Item(14, 10).Value = 2.86
$ws.Cells.Item(14, 22).Value = 1.4
$ws.Cells.Item(14, 25).Value = 11
$ws.Cells.Item(14, 27).Value = 80
$ws.Cells.Item(14, 28).Value = 9.4
$ws.Cells.Item(14, 30).Value = 18
$ws.Cells.Item(14, 33).Value = 15.5
$ws.Cells.Item(14, 36).Value = 50
$ws.Cells.Item(14, 38).Value = 85

# Row 15
$ws.Cells.Item(15, 7).Value = 2.58
$ws.Cells.Item(15, 9).Value = 4.4
$ws.Cells.Item(15, 13).Value = 1.08
$ws.Cells.Item(15, 14).Value = 2.8
$ws.Cells.Item(15, 16).Value = 1.73
$ws.Cells.Item(15, 17).Value = 1.93
$ws.Cells.Item(15, 19).Value = 3.4
$ws.Cells.Item(15, 21).Value = 1.94
$ws.Cells.Item(15, 22).Value = 1.29
$ws.Cells.Item(15, 23).Value = 1.66

# Row 16
$ws.Cells.Item(16, 7).Value = 2.94
$ws.Cells.Item(16, 9).Value = 3.7
$ws.Cells.Item(16, 10).Value = 2.74
$ws.Cells.Item(16, 14).Value = 2.9
$ws.Cells.Item(16, 16).Value = 1.65
$ws.Cells.Item(16, 17).Value = 2.08
$ws.Cells.Item(16, 19).Value = 3.8
$ws.Cells.Item(16, 20).Value = 1.88
$ws.Cells.Item(16, 22).Value = 1.37
$ws.Cells.Item(16, 23).Value = 1.51
$ws.Cells.Item(16, 24).Value = 11
$ws.Cells.Item(16, 25).Value = 11
$ws.Cells.Item(16, 37).Value = 1000
$ws.Cells.Item(16, 38).Value = 65

# Row 17
$ws.Cells.Item(17, 6).Value = 1.63
$ws.Cells.Item(17, 7).Value = 1.75
$ws.Cells.Item(17, 8).Value = 6.2
$ws.Cells.Item(17, 9).Value = 7.6
$ws.Cells.Item(17, 11).Value = 4
$ws.Cells.Item(17, 12).Value = 1.48
$ws.Cells.Item(17, 14).Value = 3
$ws.Cells.Item(17, 15).Value = 1.39
$ws.Cells.Item(17, 17).Value = 2.18
$ws.Cells.Item(17, 19).Value = 4.1
$ws.Cells.Item(17, 20).Value = 2.1
$ws.Cells.Item(17, 22).Value = 1.15
$ws.Cells.Item(17, 23).Value = 2.32
$ws.Cells.Item(17, 26).Value = 55
$ws.Cells.Item(17, 27).Value = 220
$ws.Cells.Item(17, 28).Value = 7
$ws.Cells.Item(17, 29).Value = 9
$ws.Cells.Item(17, 31).Value = 120
$ws.Cells.Item(17, 32).Value = 9
$ws.Cells.Item(17, 34).Value = 28
$ws.Cells.Item(17, 35).Value = 130
$ws.Cells.Item(17, 39).Value = 200
$ws.Cells.Item(17, 40).Value = 14
$ws.Cells.Item(17, 41).Value = 180

# Row 18
$ws.Cells.Item(18, 7).Value = 1.72
$ws.Cells.Item(18, 9).Value = 7.4
$ws.Cells.Item(18, 19).Value = 3.75
$ws.Cells.Item(18, 22).Value = 1.17
$ws.Cells.Item(18, 23).Value = 2.38
$ws.Cells.Item(18, 24).Value = 16.5
$ws.Cells.Item(18, 25).Value = 21
$ws.Cells.Item(18, 26).Value = 60
$ws.Cells.Item(18, 28).Value = 8.800000000000001
$ws.Cells.Item(18, 29).Value = 10
$ws.Cells.Item(18, 30).Value = 29
$ws.Cells.Item(18, 32).Value = 11.5
$ws.Cells.Item(18, 33).Value = 12
$ws.Cells.Item(18, 34).Value = 28
$ws.Cells.Item(18, 36).Value = 19.5
$ws.Cells.Item(18, 37).Value = 22
$ws.Cells.Item(18, 38).Value = 55
$ws.Cells.Item(18, 40).Value = 14.5

# Row 19
$ws.Cells.Item(19, 6).Value = 1.86
$ws.Cells.Item(19, 7).Value = 1.94
$ws.Cells.Item(19, 9).Value = 6.4
$ws.Cells.Item(19, 11).Value = 3.65
$ws.Cells.Item(19, 15).Value = 1.49
$ws.Cells.Item(19, 17).Value = 2.4
$ws.Cells.Item(19, 19).Value = 4.7
$ws.Cells.Item(19, 22).Value = 1.19
$ws.Cells.Item(19, 25).Value = 17.5
$ws.Cells.Item(19, 26).Value = 55
$ws.Cells.Item(19, 27).Value = 220
$ws.Cells.Item(19, 28).Value = 7.8
$ws.Cells.Item(19, 29).Value = 9.6
$ws.Cells.Item(19, 30).Value = 27
$ws.Cells.Item(19, 31).Value = 130
$ws.Cells.Item(19, 32).Value = 12
$ws.Cells.Item(19, 33).Value = 13
$ws.Cells.Item(19, 34).Value = 32
$ws.Cells.Item(19, 35).Value = 140
$ws.Cells.Item(19, 36).Value = 23
$ws.Cells.Item(19, 37).Value = 30
$ws.Cells.Item(19, 39).Value = 250
$ws.Cells.Item(19, 40).Value = 24

# Row 20
$ws.Cells.Item(20, 8).Value = 21
$ws.Cells.Item(20, 14).Value = 4.5
$ws.Cells.Item(20, 16).Value = 2.24
$ws.Cells.Item(20, 18).Value = 1.49
$ws.Cells.Item(20, 19).Value = 2.56
$ws.Cells.Item(20, 21).Value = 1.41
